$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 with the new (rounded) measurement values.
# Column A keeps its date/time serial, but with a refined value.
$ws.Range("A5").Value = 44781.9027662037

$newValues = @{
    "B5" = 8.87
    "C5" = 6.74
    "D5" = 0.18
    "E5" = 19.66
    "F5" = 15.97
    "G5" = 6.76
    "H5" = 29.2
    "I5" = 10.9
    "J5" = 5.22
    "K5" = 7.29
    "L5" = 7.96
    "M5" = 8.44
    "N5" = 2.27
    "O5" = 7.12
    "P5" = 9.98
    "Q5" = 6.05
    "R5" = 0.09
    "S5" = 0.42
    "T5" = 101.4
    "U5" = 20.02
    "V5" = 6.57
    "W5" = 13.27
    "X5" = 7.3
    "Y5" = 1.08
    "Z5" = 14.22
    "AA5" = 5.75
    "AB5" = 5.5
    "AC5" = 6.52
    "AD5" = 8.57
    "AE5" = 0.33
    "AF5" = 26.5
    "AG5" = 3.94
    "AH5" = 8.07
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# Remove row 6 entirely (it is no longer part of the dataset).
$ws.Rows.Item(6).Delete()
